{"js": "// Update the worksheet date header and every \"A\u00d7B=\" multiplication\n// prompt in the table to the new values from the regenerated output.\n// Every old value below is unique within the document, so a simple\n// exact-text search + replace for each pair is unambiguous.\nconst replacements = [\n  [\"2024-09-27 Friday\", \"2024-09-28 Saturday\"],\n  [\"813\u00d77=\", \"918\u00d77=\"],\n  [\"161\u00d77=\", \"469\u00d74=\"],\n  [\"558\u00d79=\", \"184\u00d74=\"],\n  [\"595\u00d78=\", \"330\u00d78=\"],\n  [\"283\u00d73=\", \"666\u00d72=\"],\n  [\"102\u00d75=\", \"238\u00d72=\"],\n  [\"579\u00d79=\", \"824\u00d79=\"],\n  [\"455\u00d75=\", \"232\u00d73=\"],\n  [\"819\u00d73=\", \"575\u00d79=\"],\n  [\"784\u00d79=\", \"205\u00d78=\"],\n  [\"676\u00d75=\", \"524\u00d78=\"],\n  [\"480\u00d79=\", \"473\u00d79=\"],\n  [\"206\u00d72=\", \"742\u00d74=\"],\n  [\"405\u00d75=\", \"360\u00d75=\"],\n  [\"549\u00d73=\", \"683\u00d78=\"],\n  [\"426\u00d75=\", \"751\u00d76=\"],\n  [\"805\u00d74=\", \"369\u00d76=\"],\n  [\"568\u00d73=\", \"303\u00d78=\"],\n  [\"536\u00d72=\", \"423\u00d75=\"],\n  [\"874\u00d74=\", \"454\u00d76=\"],\n  [\"554\u00d72=\", \"924\u00d72=\"],\n  [\"344\u00d73=\", \"469\u00d79=\"],\n  [\"405\u00d73=\", \"399\u00d72=\"],\n  [\"976\u00d77=\", \"913\u00d73=\"],\n  [\"112\u00d74=\", \"987\u00d77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date header and every \"A\u00d7B=\" multiplication\n# prompt in the table to the new values from the regenerated output.\n# Every old value is unique within the document, so Find/Replace on\n# each exact pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-09-27 Friday\", \"2024-09-28 Saturday\"),\n  @(\"813\u00d77=\", \"918\u00d77=\"),\n  @(\"161\u00d77=\", \"469\u00d74=\"),\n  @(\"558\u00d79=\", \"184\u00d74=\"),\n  @(\"595\u00d78=\", \"330\u00d78=\"),\n  @(\"283\u00d73=\", \"666\u00d72=\"),\n  @(\"102\u00d75=\", \"238\u00d72=\"),\n  @(\"579\u00d79=\", \"824\u00d79=\"),\n  @(\"455\u00d75=\", \"232\u00d73=\"),\n  @(\"819\u00d73=\", \"575\u00d79=\"),\n  @(\"784\u00d79=\", \"205\u00d78=\"),\n  @(\"676\u00d75=\", \"524\u00d78=\"),\n  @(\"480\u00d79=\", \"473\u00d79=\"),\n  @(\"206\u00d72=\", \"742\u00d74=\"),\n  @(\"405\u00d75=\", \"360\u00d75=\"),\n  @(\"549\u00d73=\", \"683\u00d78=\"),\n  @(\"426\u00d75=\", \"751\u00d76=\"),\n  @(\"805\u00d74=\", \"369\u00d76=\"),\n  @(\"568\u00d73=\", \"303\u00d78=\"),\n  @(\"536\u00d72=\", \"423\u00d75=\"),\n  @(\"874\u00d74=\", \"454\u00d76=\"),\n  @(\"554\u00d72=\", \"924\u00d72=\"),\n  @(\"344\u00d73=\", \"469\u00d79=\"),\n  @(\"405\u00d73=\", \"399\u00d72=\"),\n  @(\"976\u00d77=\", \"913\u00d73=\"),\n  @(\"112\u00d74=\", \"987\u00d77=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
